# Insert a new "Match ID" column at the very left of the data table.
# This shifts every existing column (A:AC) one place to the right (B:AD)
# and populates the new column A with:
#   - row 3 (header row): "Match ID"
#   - rows 4-20 (data rows, including the hidden totals row 20): 21
# The header cell and the per-player rows (3-19) pick up a bold,
# border-less style (matches the new cellXfs entry); the hidden summary
# row (20) is left with the default style, matching how Excel only
# copies the format of the row immediately above when a new column is
# inserted through the grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything right by inserting a fresh column before column A.
$ws.Columns.Item(1).Insert() | Out-Null

# New header + data values in column A.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A4:A20").Value = 21

# Bold the header + visible player rows (not the hidden totals row).
$ws.Range("A3:A19").Font.Bold = $true

# Re-fit the hidden summary row so it doesn't pick up a stray explicit
# row height from the value write (keeps row 20 identical to before).
$ws.Rows.Item(20).AutoFit()

# Mirror the author's final selection: the whole new column's data span.
$ws.Range("A3:A19").Select() | Out-Null
